# Update crypto price/volume data per the daily GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose numeric-looking values have a significant
# trailing zero that plain numeric entry would otherwise drop (e.g. "2.660").
$ws.Range("D38,D44,D46,D49").NumberFormat = "@"

# Price (D) and Volume(1h) (E) updates, row by row.
$ws.Range("D2").Value = "30.721.40"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.900.64"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "239.18"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4816"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "0.06557"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "1.915.45"
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("D11").Value = "0.07466"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "16.74"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "5.123"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "0.6672"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "30.680.19"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "0.000007632"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "231.89"
$ws.Range("E20").Value = "  +6.51%  "
$ws.Range("D21").Value = "2.150.04"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").Value = "5.308"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "6.244"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "170.25"
$ws.Range("E25").Value = "  +3.95%  "
$ws.Range("D26").Value = "9.322"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").Value = "18.72"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "1.969"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  +9.37%  "
$ws.Range("D31").Value = "4.366"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "0.05093"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "1.219"
$ws.Range("E34").Value = "  +7.21%  "
$ws.Range("D35").Value = "0.7582"
$ws.Range("D36").Value = "2.714"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "0.01881"
$ws.Range("E37").Value = "  +2.88%  "
$ws.Range("D38").Value = "2.660"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").Value = "0.9191"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "2.089"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").Value = "107.04"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "0.4307"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "5.770"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").Value = "7.445"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "64.40"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "0.1276"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").Value = "1.489"
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("D49").Value = "9.020"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "0.05673"
$ws.Range("E51").Value = "  -0.05%  "
